$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 2 (alpha_distance_range)
$ws.Range("B2").Value = 5

# Update values in row 3 (beta_distance_range)
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 9

# Update values in row 4 (ratio_threshold_range)
$ws.Range("C4").Value = 1.4

# Row 5 (previously theta_threshold_range) now becomes pie_threshold_range with new values
$ws.Range("A5").Value = "pie_threshold_range"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

# C5 previously used a special Times New Roman style; reset it to match the other
# data cells (same formatting as B5, which already uses the common data style).
$ws.Range("B5").Copy()
$ws.Range("C5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C5").Value = 15
$excel.CutCopyMode = $false

# Remove the now-duplicate last row (previously pie_threshold_range)
$ws.Rows.Item(6).Delete()

$ws.Range("C5").Select()
